$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Folds" column (B) values by 1 for rows 2-6
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 5
